$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.742.72'
$ws.Range('E2').Value = '  +0.49%  '

$ws.Range('D3').Value = '3.419.29'
$ws.Range('E3').Value = '  +0.83%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '407.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.88%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '127.92'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.78%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.629'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +7.00%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.731'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.61%  '

$ws.Range('E10').Value = '  +8.82%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.47'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.48%  '

$ws.Range('E12').Value = '  +9.59%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.953.54'
$ws.Range('E14').Value = '  +0.73%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +8.07%  '

$ws.Range('E16').Value = '  +41.12%  '

$ws.Range('D17').Value = '3.390.77'
$ws.Range('E17').Value = '  +0.15%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.20'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.15%  '

$ws.Range('D20').Value = '61.710.83'
$ws.Range('E20').Value = '  +0.46%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '442.76'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +42.74%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '91.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +10.15%  '

$ws.Range('E23').Value = '  +0.99%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.92'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.11%  '

$ws.Range('E25').Value = '  +2.77%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '32.79'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.81%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.75%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.45%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.65'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.21%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.72'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.68%  '

$ws.Range('E31').Value = '  +6.07%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.170'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.42%  '

$ws.Range('E33').Value = '  +0.39%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '42.43'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.47%  '

$ws.Range('E35').Value = '  -0.08%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0496'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.50%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.18'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.88%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.09%  '

$ws.Range('E39').Value = '  +0.62%  '

$ws.Range('E40').Value = '  +8.06%  '

$ws.Range('E41').Value = '  -0.78%  '

$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.311'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.17%  '

$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '141.98'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.59%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.24%  '

$ws.Range('E45').Value = '  +1.59%  '

$ws.Range('E46').Value = '  +13.74%  '

$ws.Range('E47').Value = '  -0.61%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.18'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.38%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.15'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +14.04%  '

$ws.Range('D50').Value = '3.766.45'
$ws.Range('E50').Value = '  +0.94%  '

$ws.Range('D51').Value = '2.117.83'
$ws.Range('E51').Value = '  +1.41%  '
